# Rimosse domande/dubbi: elimina i paragrafi finali (la riga "..." e le
# tre domande/dubbi che la seguono, oltre al paragrafo vuoto finale)
# lasciando come ultimo paragrafo del corpo quello che termina con
# "Formulare e modellare il problema in forma discreta."

$d = $word.ActiveDocument

$anchorText = "Formulare e modellare il problema in forma discreta."

$find = $d.Content.Find
$found = $find.Execute($anchorText, $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)

if ($found) {
    $anchorParagraph = $find.Parent
    # Espande la selezione fino a includere l'intero paragrafo (compreso
    # il relativo segno di paragrafo) così il punto di taglio è esatto.
    $anchorParagraph.Expand(4) | Out-Null

    $deleteStart = $anchorParagraph.End
    $deleteEnd = $d.Content.End

    if ($deleteEnd -gt $deleteStart) {
        $trailingRange = $d.Range($deleteStart, $deleteEnd)
        $trailingRange.Delete()
    }
}
